$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2 (AccountDetails) — fill in the guest / register-user test values
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = "avayugundla@helenoftroy.com"
$ws.Range("C2").Value = "avayugundla@helenoftroy.com"
$ws.Range("D2").Value = "Lotuswave@123"
$ws.Range("E2").Value = "Lotuswave@123"
$ws.Range("F2").Value = "Test"
$ws.Range("G2").Value = "Qa"
$ws.Range("K2").Value = "844 N Colony Rd"
$ws.Range("L2").Value = "Wallingford"
$ws.Range("M2").Value = "United States"
$ws.Range("N2").Value = "Connecticut"
$ws.Range("O2").Value = "'06492"
$ws.Range("P2").Value = 9898989898

# ---------------------------------------------------------------------------
# Row 17 (new "Address" test-data row, mirrors row 16's "Colorado Address")
# ---------------------------------------------------------------------------
$ws.Range("A17").Value = "Address"
$ws.Range("F17").Value = "Qa"
$ws.Range("G17").Value = "Test"
$ws.Range("H17").Value = "avayugundla@helenoftroy.com"
$ws.Range("K17").Value = "844 N Colony Rd"
$ws.Range("L17").Value = "Wallingford"
$ws.Range("M17").Value = "United States"
$ws.Range("N17").Value = "Connecticut"
$ws.Range("O17").Value = "'06492"
$ws.Range("P17").Value = 9898989898

# ---------------------------------------------------------------------------
# Hyperlinks — e-mail / password fields auto-linked (mailto:), new address
# row's e-mail hyperlinked like the existing row 16 one
# ---------------------------------------------------------------------------
$null = $ws.Hyperlinks.Add($ws.Range("B2"), "mailto:avayugundla@helenoftroy.com")
$null = $ws.Hyperlinks.Add($ws.Range("C2"), "mailto:avayugundla@helenoftroy.com")
$null = $ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Lotuswave@123")
$null = $ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Lotuswave@123")
$null = $ws.Hyperlinks.Add($ws.Range("H17"), "mailto:avayugundla@helenoftroy.com")

# ---------------------------------------------------------------------------
# Selection moves to C4 (matches the saved sheetView selection in the diff)
# ---------------------------------------------------------------------------
$null = $ws.Range("C4").Select()
